$d = $word.ActiveDocument

$d.Content.Find.Execute("external_short_name", $true, $false, $false, $false, $false, $true, 1, $false, "venue_name", 2)
